$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two new columns -------------------------------------------------
# New "gas" column lands at D (pushes former D..AZ right by one).
$ws.Columns("D").Insert()
# New "measurement" column lands at I (former H becomes I after the shift above).
$ws.Columns("I").Insert()

# --- Insert a new row under row 2 for the second gas measurement ----------------
# (row 2 = CH4 reading; the new row 3 duplicates it for the CO2 reading)
$ws.Rows("3").Insert()

# Copy row 2's existing values (A,B,C,E,F,J,K,L,M,N) down into the newly
# inserted row 3, explicitly (avoids COM clipboard float round-off).
$ws.Range("A3").Value2 = $ws.Range("A2").Value2
$ws.Range("B3").Value2 = $ws.Range("B2").Value2
$ws.Range("C3").Value2 = $ws.Range("C2").Value2
$ws.Range("E3").Value2 = $ws.Range("E2").Value2
$ws.Range("F3").Value2 = $ws.Range("F2").Value2
$ws.Range("J3").Value2 = $ws.Range("J2").Value2
$ws.Range("K3").Value2 = $ws.Range("K2").Value2
$ws.Range("L3").Value2 = $ws.Range("L2").Value2
$ws.Range("M3").Value2 = $ws.Range("M2").Value2
$ws.Range("N3").Value2 = $ws.Range("N2").Value2

# --- Header label + values for the "gas" column (order matches shared-string
#     insertion order: gas, CH4, CO2, measurement) -------------------------------
$ws.Range("D1").Value2 = "gas"

# --- Fill in the "gas" reading type for every data row (alternating CH4/CO2) ----
$ws.Range("D2").Value2 = "CH4"
$ws.Range("D3").Value2 = "CO2"
$ws.Range("D4").Value2 = "CH4"
$ws.Range("D5").Value2 = "CO2"
$ws.Range("D6").Value2 = "CH4"
$ws.Range("D7").Value2 = "CO2"

# --- Header label for the "measurement" column (added last) --------------------
$ws.Range("I1").Value2 = "measurement"

# --- Match the saved selection / active cell recorded in the edit ---------------
$ws.Range("D8").Select() | Out-Null
